$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the weekly price-report data between row 2 and row 3:
# row 2 becomes what row 3 used to hold (and vice versa) for the
# date, volume, min/max/avg price and $/Kg columns.

$ws.Range("D2").Value = 44804
$ws.Range("J2").Value = 50
$ws.Range("K2").Value = 9500
$ws.Range("M2").Value = 9750
$ws.Range("P2").Value = 542

$ws.Range("D3").Value = 44714
$ws.Range("J3").Value = 80
$ws.Range("K3").Value = 9000
$ws.Range("M3").Value = 9500
$ws.Range("P3").Value = 528
